$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): re-write as literal numbers 0..9 (bug introduced by the
#     exporter while fixing the PPG column issue - header cells lose their
#     text and end up holding their column index instead) ---
for ($i = 0; $i -le 9; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# --- Row 2 / Row 3: rotate the nome / cv_lattes / PPG columns (A,B,C) ---
# Before:  A=nome             B=cv_lattes                               C=PPG
# After:   A=cv_lattes        B=PPG                                     C=nome
$ws.Range("A2").Value = "http://lattes.cnpq.br/3269118444404338"
$ws.Range("B2").Value = "GEAS"
$ws.Range("C2").Value = "Patricia Storópoli Tzortzis"

$ws.Range("A3").Value = "http://lattes.cnpq.br/2281909649311607"
$ws.Range("B3").Value = "PPGA-D"
$ws.Range("C3").Value = "José Eduardo Storopoli"

# --- Hyperlinks move from B2/B3 to A2/A3. The export still links both rows
#     to the same (first) cv_lattes URL, exactly like before the fix. ---
$ws.Range("A1:J3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "http://lattes.cnpq.br/3269118444404338")
$ws.Hyperlinks.Add($ws.Range("A3"), "http://lattes.cnpq.br/3269118444404338")

# --- Styles: hyperlink style follows the cv_lattes column (now A), and the
#     old B column reverts to the default (Normal) style ---
$ws.Range("A2:A3").Style = "Hyperlink"
$ws.Range("B2:B3").Style = "Normal"
